$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 58: 10/29/2025 profit allocation data.
# Force column A to be treated as text (like the existing date-string
# cells above it) instead of letting Excel auto-convert "10/29/2025"
# into a date serial number.
$ws.Cells.Item(58, 1).NumberFormat = "@"
$ws.Cells.Item(58, 1).Value = "10/29/2025"
$ws.Cells.Item(58, 1).Style = "Normal"

$ws.Cells.Item(58, 2).Value = 0.1825660000012012
$ws.Cells.Item(58, 3).Value = 0.8174339999987988
